$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table (rows 3-5) ---
# Row 3: new FastConnect 6900 ...1229 entry inserted at the top
$ws.Range("A3").Value = "Qualcomm FastConnect 6900 Wi-Fi 6E Dual Band Simultaneous (DBS) WiFiCx Network Adapter - 2.0.0.1229"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 88.90000000000001

# Row 4: previously row 3's driver, with refreshed counts
$ws.Range("A4").Value = "Qualcomm Atheros AR9580 Wireless Network Adapter - 10.1.10.5"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 878
$ws.Range("D4").Value = 94.8

# Row 5: previously row 4's driver, with refreshed counts (Realtek row is gone)
$ws.Range("A5").Value = "Qualcomm Atheros AR9580 Wireless Network Adapter - 3.0.2.201"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1058
$ws.Range("D5").Value = 95.3

# Totals row
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 1937

# --- Good Drivers table (rows 14-15, new entries) ---
# Row 14
$ws.Range("A14").Value = "Qualcomm FastConnect 6900 Wi-Fi 6E Dual Band Simultaneous (DBS) WiFiCx Network Adapter - 2.0.0.1277"
$ws.Range("B14").Value = 28693
$ws.Range("B14").HorizontalAlignment = -4152
$ws.Range("B14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("D14").HorizontalAlignment = -4152
$ws.Range("E14").HorizontalAlignment = -4152
$ws.Range("E14").Formula = "=""2024-11-18"""
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4163) | Out-Null

# Row 15
$ws.Range("A15").Value = "Qualcomm FastConnect 6900 Wi-Fi 6E Dual Band Simultaneous (DBS) WiFiCx Network Adapter - 2.0.0.1229"
$ws.Range("B15").Value = 196400
$ws.Range("B15").HorizontalAlignment = -4152
$ws.Range("B15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("D15").HorizontalAlignment = -4152
$ws.Range("E15").HorizontalAlignment = -4152
$ws.Range("E15").Formula = "=""2024-08-26"""
$ws.Range("E15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# --- Column A widened to fit the longer driver names ---
$ws.Columns.Item(1).ColumnWidth = 100.14
